# Insert a new weekly data row for "Vega Modelo de Temuco" / Alcachofa.
# The new record is inserted at row 257, pushing all the existing rows
# from 257..324 down to 258..325 (dimension grows from A1:R324 to A1:R325).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 257..324 down by one to make room for the new record.
$ws.Rows("257:257").Insert()

# Populate the newly inserted row 257 with the new weekly record.
$ws.Range("A257").Value = 10
$ws.Range("B257").Value = "Vega Modelo de Temuco"
$ws.Range("C257").Value = "La Araucanía"
$ws.Range("D257").Value = 45135
$ws.Range("E257").Value = 9
$ws.Range("F257").Value = 100112013
$ws.Range("G257").Value = "Alcachofa"
$ws.Range("H257").Value = "Madrigal"
$ws.Range("I257").Value = "Primera"
$ws.Range("J257").Value = 55
$ws.Range("K257").Value = 14000
$ws.Range("L257").Value = 14000
$ws.Range("M257").Value = 14000
$ws.Range("N257").Value = "$/caja 40 unidades"
$ws.Range("O257").Value = "Provincia de Limarí"
$ws.Range("P257").Value = 350
$ws.Range("Q257").Value = 40
$ws.Range("R257").Value = "Hortaliza"
